$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 116
$ws1.Range("F3").Value = 7582
$ws1.Range("F7").Value = 4215
$ws1.Range("F8").Value = 327
$ws1.Range("F9").Value = 586
$ws1.Range("F12").Value = 156

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 116
$ws4.Range("F4").Value = 7582
$ws4.Range("F9").Value = 4215
$ws4.Range("F10").Value = 327
$ws4.Range("F11").Value = 586
$ws4.Range("F15").Value = 156
